# Generate Report for Handoff
#
# Replaces the placeholder source-document GUID
# (80575781-ce0e-481f-ae96-719873221745) with the freshly generated one
# (5227ef20-1367-4ec0-9adc-8f48080d7d0e), refreshes the xliff hash tokens,
# and bumps the handoff timestamps on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldGuidHash = "80575781-ce0e-481f-ae96-719873221745"
$newGuidHash = "5227ef20-1367-4ec0-9adc-8f48080d7d0e"

$oldXliffHash = "a0d56a9363ab0fb6678645ac25ad51ac03a3c016"
$newXliffHash = "0c1dc2013969902c6a7fa565b0150faafc4b5d80"

# -- File Name / Path And Name (Overview, zh-cn, de-de : column A) ----------
$newFileName = "$newGuidHash.md"
$newPathAndName = "e2e\$newGuidHash.md"

$wsOverview.Range("A2").Value = $newFileName
$wsZhCn.Range("A2").Value     = $newFileName
$wsDeDe.Range("A2").Value     = $newFileName

# Path And Name + its hyperlink display text (Overview sheet, column B)
$wsOverview.Range("B2").Value = $newPathAndName

# -- Latest Handoff File (per-language "Latest Handoff File" column G) -----
$wsZhCn.Range("G2").Value = "$newGuidHash.$newXliffHash.zh-cn.xlf"
$wsDeDe.Range("G2").Value = "$newGuidHash.$newXliffHash.de-de.xlf"

# -- Latest Handoff Datetime (per-language column H) ------------------------
$wsZhCn.Range("H2").Value = "2016-08-24 15:07:15"
$wsDeDe.Range("H2").Value = "2016-08-24 15:07:20"

# -- Latest HO Xliff Generate Date (Overview column G) ----------------------
# This mirrors the most recent per-language handoff datetime (de-de, the
# last one generated).
$wsOverview.Range("G2").Value = "2016-08-24 15:07:20"

# -- Update hyperlink display text so it matches the new file names --------
$overviewLinks = @($wsOverview.Hyperlinks)
$overviewLinks[0].TextToDisplay = $newPathAndName

$zhCnLinks = @($wsZhCn.Hyperlinks)
$zhCnLinks[0].TextToDisplay = $newFileName

$deDeLinks = @($wsDeDe.Hyperlinks)
$deDeLinks[0].TextToDisplay = $newFileName
